# Update the benchmark table for the first draft:
#  - split each "<Period> Base" column into "<Period> Base mean" / "<Period> Base std"
#  - rename algorithm "CART" -> "DTREE"
#  - drop the "NB" algorithm row entirely
#  - refresh every numeric value with the new run's results

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the last data row (row 9, algorithm "NB") before anything else so
#    subsequent row numbers line up with the final 8-row (1 header + 7 data)
#    layout.
# ---------------------------------------------------------------------------
$ws.Rows("9").Delete()

# ---------------------------------------------------------------------------
# 2) Extend the header row with the five new "std" columns (H1:L1), copying
#    the existing header formatting (bold, border, centered) from G1.
# ---------------------------------------------------------------------------
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Header labels (B1:L1).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# ---------------------------------------------------------------------------
# 4) Algorithm names in column B (row 5's "CART" becomes "DTREE"; the old
#    "NB" row is already gone).
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "LR"
$ws.Range("B3").Value = "LDA"
$ws.Range("B4").Value = "KNN"
$ws.Range("B5").Value = "DTREE"
$ws.Range("B6").Value = "RTREE"
$ws.Range("B7").Value = "XTREE"
$ws.Range("B8").Value = "SVM"

# ---------------------------------------------------------------------------
# 5) Refreshed numeric results, columns C:L, rows 2:8.
# ---------------------------------------------------------------------------
$data = @(
    @(0.8227043832088959, 0.01862514451755787, 0.81680997916091,   0.02016003136401769, 0.8141002689724395, 0.02698399373550396, 0.8087085764982754, 0.02330693759478715, 0.791492623450966,  0.01627445701162525),
    @(0.8278176516051394, 0.01716060341196075, 0.8145849899016977, 0.02228374804388591, 0.8071371820652065, 0.01873490343086296, 0.7971884962727798, 0.02979905705404396, 0.7908893920597986, 0.01946116434771431),
    @(0.7797237575740318, 0.01034381290272372, 0.7838810119296056, 0.01924338605229094, 0.7849327771403607, 0.0216377810033402,  0.7902076932162205, 0.02194352585381001, 0.780342299770113,  0.01546790320066242),
    @(0.789744216493433,  0.02177300219927225, 0.7656482131248989, 0.02152487416445197, 0.7601505968872101, 0.01628677003725389, 0.7621046429825049, 0.02809523536527091, 0.7671530738391767, 0.04684496958647222),
    @(0.8205546666999201, 0.01445369909918223, 0.8017108233902691, 0.01251135254418705, 0.7828549989833528, 0.01697743647061012, 0.7621020619911592, 0.02440502254272588, 0.7360998808234999, 0.02524221614875536),
    @(0.834209447315047,  0.01417499778868238, 0.8103393824228515, 0.02240009067659626, 0.7947514265648719, 0.01743246367610575, 0.7835512160843862, 0.03809256859415011, 0.7781105147823051, 0.0298409006013912),
    @(0.8353244815112054, 0.01430710421391247, 0.8247116312815799, 0.01744686961647583, 0.8167438931598692, 0.02061544725623895, 0.8118694813104035, 0.0248090676157972,  0.7967845988153696, 0.01910166225212886)
)

$startCol = 3   # column C
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $startCol + $j).Value = $values[$j]
    }
}
